# Auto-generated edit script: updates numeric cells per the target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 32850
$ws.Range("J3").Value = 32850
$ws.Range("L3").Value = 32850
$ws.Range("N3").Value = -33078

$ws.Range("H15").Value = 347.69
$ws.Range("I15").Value = 347.69
$ws.Range("K15").Value = 1043.07
$ws.Range("M15").Value = -874.0699999999999

$ws.Range("H28").Value = 353.05884
$ws.Range("I28").Value = 353.05884
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 353.05884
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 131.94116
$ws.Range("N28").ClearContents()

$ws.Range("H62").Value = 10324.296
$ws.Range("I62").Value = 13079.765
$ws.Range("J62").Value = 5640
$ws.Range("K62").Value = 13079.765
$ws.Range("L62").Value = 5640
$ws.Range("M62").Value = -12455.765
$ws.Range("N62").Value = -6888

$ws.Range("H65").Value = 10324.296
$ws.Range("I65").Value = 13079.765
$ws.Range("J65").Value = 5640
$ws.Range("K65").Value = 65398.825
$ws.Range("L65").Value = 28200
$ws.Range("M65").Value = -62278.825
$ws.Range("N65").Value = -34440

$ws.Range("H102").Value = 32850
$ws.Range("J102").Value = 32850
$ws.Range("L102").Value = 32850
$ws.Range("N102").Value = -39340

$ws.Range("H105").Value = 34653.332
$ws.Range("J105").Value = 34653.332
$ws.Range("L105").Value = 34653.332
$ws.Range("N105").Value = -41641.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1519622.6
$ws.Range("I32").Value = 1856559.1
$ws.Range("J32").Value = 3408.3
$ws.Range("K32").Value = 1856559.1
$ws.Range("L32").Value = 3408.3
$ws.Range("M32").Value = -1856272.1
$ws.Range("N32").Value = -3982.3

$ws.Range("H92").Value = 34333
$ws.Range("J92").Value = 34333
$ws.Range("L92").Value = 34333
$ws.Range("N92").Value = -39325

$ws.Range("H119").Value = 34990
$ws.Range("J119").Value = 34990
$ws.Range("L119").Value = 34990
$ws.Range("N119").Value = -44666

$ws.Range("H124").Value = 17714.5
$ws.Range("J124").Value = 17714.5
$ws.Range("L124").Value = 17714.5
$ws.Range("N124").Value = -27534.5

$ws.Range("H125").Value = 33667.5
$ws.Range("J125").Value = 33667.5
$ws.Range("L125").Value = 33667.5
$ws.Range("N125").Value = -43507.5

$ws.Range("H132").Value = 35519.062
$ws.Range("I132").Value = 96937.45
$ws.Range("J132").Value = 3347.524
$ws.Range("K132").Value = 290812.35
$ws.Range("L132").Value = 10042.572
$ws.Range("M132").Value = -288282.35
$ws.Range("N132").Value = -15102.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 42930
$ws.Range("J130").Value = 42930
$ws.Range("L130").Value = 42930
$ws.Range("N130").Value = -52970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2746.7812
$ws.Range("I31").Value = 1258.2084
$ws.Range("J31").Value = 7212.5
$ws.Range("K31").Value = 1258.2084
$ws.Range("L31").Value = 7212.5
$ws.Range("M31").Value = -963.2084
$ws.Range("N31").Value = -7802.5

$ws.Range("H34").Value = 2746.7812
$ws.Range("I34").Value = 1258.2084
$ws.Range("J34").Value = 7212.5
$ws.Range("K34").Value = 1258.2084
$ws.Range("L34").Value = 7212.5
$ws.Range("M34").Value = -1056.2084
$ws.Range("N34").Value = -7616.5

$ws.Range("H58").Value = 4188.425
$ws.Range("I58").Value = 6116.316
$ws.Range("J58").Value = 2444.1428
$ws.Range("K58").Value = 6116.316
$ws.Range("L58").Value = 2444.1428
$ws.Range("M58").Value = -5913.316
$ws.Range("N58").Value = -2850.1428

$ws.Range("H107").Value = 194.5
$ws.Range("I107").Value = 106.3125
$ws.Range("K107").Value = 106.3125
$ws.Range("M107").Value = 1813.6875

$ws.Range("H132").Value = 2811.7
$ws.Range("I132").Value = 1081.2
$ws.Range("J132").Value = 4542.2
$ws.Range("K132").Value = 3243.6
$ws.Range("L132").Value = 13626.6
$ws.Range("M132").Value = -713.6000000000004
$ws.Range("N132").Value = -18686.6

$ws.Range("H136").Value = 4188.425
$ws.Range("I136").Value = 6116.316
$ws.Range("J136").Value = 2444.1428
$ws.Range("K136").Value = 18348.948
$ws.Range("L136").Value = 7332.428400000001
$ws.Range("M136").Value = -15798.948
$ws.Range("N136").Value = -12432.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71.26667
$ws.Range("I2").Value = 28.333334
$ws.Range("J2").Value = 99.888885
$ws.Range("K2").Value = 170.000004
$ws.Range("L2").Value = 599.33331
$ws.Range("M2").Value = -57.00000399999999
$ws.Range("N2").Value = -825.33331

$ws.Range("H38").Value = 99.78570999999999
$ws.Range("I38").Value = 65
$ws.Range("J38").Value = 162.4
$ws.Range("K38").Value = 195
$ws.Range("L38").Value = 487.2
$ws.Range("M38").Value = 152
$ws.Range("N38").Value = -1181.2

$ws.Range("H131").Value = 2437.2104
$ws.Range("I131").Value = 4770
$ws.Range("J131").Value = 1999.8125
$ws.Range("K131").Value = 14310
$ws.Range("L131").Value = 5999.4375
$ws.Range("M131").Value = -9270
$ws.Range("N131").Value = -16079.4375

$ws.Range("H134").Value = 2988.0435
$ws.Range("I134").Value = 2955.682
$ws.Range("J134").Value = 3700
$ws.Range("K134").Value = 8867.045999999998
$ws.Range("L134").Value = 11100
$ws.Range("M134").Value = -3797.045999999998
$ws.Range("N134").Value = -21240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 4950
$ws.Range("J39").Value = 4950
$ws.Range("L39").Value = 4950
$ws.Range("N39").Value = -6014

$ws.Range("H132").Value = 5090.3335
$ws.Range("I132").Value = 5660.8335
$ws.Range("J132").Value = 4519.8335
$ws.Range("K132").Value = 16982.5005
$ws.Range("L132").Value = 13559.5005
$ws.Range("M132").Value = -14452.5005
$ws.Range("N132").Value = -18619.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 29000
$ws.Range("J101").Value = 29000
$ws.Range("L101").Value = 29000
$ws.Range("N101").Value = -35490

$ws.Range("H104").Value = 17685
$ws.Range("J104").Value = 17685
$ws.Range("L104").Value = 17685
$ws.Range("N104").Value = -24673

$ws.Range("H132").Value = 20838054
$ws.Range("I132").Value = 50004560
$ws.Range("J132").Value = 4835.5
$ws.Range("K132").Value = 150013680
$ws.Range("L132").Value = 14506.5
$ws.Range("M132").Value = -150011150
$ws.Range("N132").Value = -19566.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 24939.25
$ws.Range("J104").Value = 24939.25
$ws.Range("L104").Value = 24939.25
$ws.Range("N104").Value = -31927.25

$ws.Range("H132").Value = 71432180
$ws.Range("I132").Value = 166669330
$ws.Range("J132").Value = 4312.25
$ws.Range("K132").Value = 500007990
$ws.Range("L132").Value = 12936.75
$ws.Range("M132").Value = -500005460
$ws.Range("N132").Value = -17996.75
